$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update repaymentstrategy value (row 17) with the new scenario text
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Activate the sheet and move the selection to B17 to match the saved view state
$ws.Activate()
$ws.Range("B17").Select()
